$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "88.886.67"
$ws.Range("E2").Value = "  -3.60%  "

# Row 3
$ws.Range("D3").Value = "3.135.54"
$ws.Range("E3").Value = "  -4.28%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.06"
$ws.Range("E5").Value = "  -0.94%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "635.54"
$ws.Range("E6").Value = "  +1.01%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.395"
$ws.Range("E7").Value = "  -4.60%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.759"
$ws.Range("E8").Value = "  +5.04%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  +0.05%  "

# Row 10
$ws.Range("D10").Value = "3.132.04"
$ws.Range("E10").Value = "  -4.20%  "

# Row 11
$ws.Range("E11").Value = "  -5.32%  "

# Row 12
$ws.Range("E12").Value = "  -0.78%  "

# Row 13
$ws.Range("E13").Value = "  -5.76%  "

# Row 14
$ws.Range("E14").Value = "  -0.58%  "

# Row 15
$ws.Range("D15").Value = "88.698.21"
$ws.Range("E15").Value = "  -3.64%  "

# Row 16
$ws.Range("D16").Value = "3.708.97"
$ws.Range("E16").Value = "  -4.42%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "32.38"
$ws.Range("E17").Value = "  -5.59%  "

# Row 18
$ws.Range("D18").Value = "3.142.68"
$ws.Range("E18").Value = "  -3.59%  "

# Row 19
$ws.Range("B19").Value = "PEPE"
$ws.Range("C19").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000230"
$ws.Range("E19").Value = "  +18.96%  "

# Row 20
$ws.Range("B20").Value = "SuiNetwork"
$ws.Range("C20").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.39"
$ws.Range("E20").Value = "  +1.64%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.29"
$ws.Range("E21").Value = "  -5.62%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "426.70"
$ws.Range("E22").Value = "  -2.95%  "

# Row 23
$ws.Range("E23").Value = "  -6.00%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.90"
$ws.Range("E24").Value = "  -7.23%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.43"
$ws.Range("E25").Value = "  +1.05%  "

# Row 26
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "82.51"
$ws.Range("E26").Value = "  +6.98%  "

# Row 27
$ws.Range("B27").Value = "Aptos"
$ws.Range("C27").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.54"
$ws.Range("E27").Value = "  -6.54%  "

# Row 28
$ws.Range("D28").Value = "3.301.75"
$ws.Range("E28").Value = "  -4.62%  "

# Row 29
$ws.Range("E29").Value = "  +0.07%  "

# Row 30
$ws.Range("E30").Value = "  -12.87%  "

# Row 31
$ws.Range("E31").Value = "  -0.40%  "

# Row 32
$ws.Range("E32").Value = "  +11.47%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.19"
$ws.Range("E33").Value = "  -6.65%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "509.13"
$ws.Range("E34").Value = "  -7.95%  "

# Row 35
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.147"
$ws.Range("E35").Value = "  +12.50%  "

# Row 36
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.04"
$ws.Range("E36").Value = "  -1.29%  "

# Row 37
$ws.Range("E37").Value = "  +1.49%  "

# Row 38
$ws.Range("E38").Value = "  -4.70%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "21.96"
$ws.Range("E39").Value = "  -3.10%  "

# Row 41
$ws.Range("E41").Value = "  +0.24%  "

# Row 42
$ws.Range("E42").Value = "  -0.05%  "

# Row 43
$ws.Range("E43").Value = "  -6.75%  "

# Row 44
$ws.Range("E44").Value = "  -7.54%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "145.24"
$ws.Range("E45").Value = "  -3.59%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "43.77"
$ws.Range("E46").Value = "  -3.71%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.129"
$ws.Range("E47").Value = "  -0.97%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "164.67"
$ws.Range("E48").Value = "  -8.65%  "

# Row 49
$ws.Range("E49").Value = "  -1.23%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.30"
$ws.Range("E50").Value = "  -3.43%  "

# Row 51
$ws.Range("E51").Value = "  -5.94%  "
